$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (V2:AO2) repeats the same values already present in B2:U2 -- copy
# the range so the floating point values round-trip exactly.
$ws.Range("B2:U2").Copy()
$ws.Range("V2").PasteSpecial()
$excel.CutCopyMode = $false

# Row 3 (V3:AO3) new 0/1 data
$row3 = @(0,0,0,1,0,1,0,0,1,0,0,1,0,1,0,0,0,0,0,1)
$cols = @("V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $row3[$i]
}

# Update the view: scroll so column L is the top-left visible column,
# and select AC8 as the active cell.
$ws.Application.ActiveWindow.ScrollColumn = 12
$ws.Range("AC8").Select()
